$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append below the existing data (row 229 is the last
# existing row). Copy the last row's formatting down first so the new date
# cells (column A) inherit the same style (date number format + alignment)
# already used by the rest of column A, then overwrite the values.
$ws.Range("A229:D229").Copy($ws.Range("A230:D233"))

$data = @(
    @(44304, 35, 260, 137.556675995831),
    @(44305, 59, 253, 133.8532270267125),
    @(44306, 17, 260, 137.556675995831),
    @(44307, 11, 259, 137.0276118573855)
)

$startRow = 230
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
